# Add a "carbon debt" factor calculation block to the biofuel (biomass
# product) production process sheets: "Charcoal" and "Syngas".
#
# Each sheet gets two new rows appended directly below the existing
# meta-process calculation rows, describing:
#   row N+1:  carbon debt | feedstock | inflows | embodied CO2 | temp | Ratio | CO2 removal
#   row N+2:              | embodied CO2 | temp | debt CO2 | outflows | Ratio | carbon debt factor

$wb = $excel.ActiveWorkbook

function Add-CarbonDebtRows {
    param($sheetName)

    $ws = $wb.Worksheets.Item($sheetName)
    $used = $ws.UsedRange
    $lastRow = $used.Row + $used.Rows.Count - 1
    if ($lastRow -lt 1) { $lastRow = 1 }

    $r1 = $lastRow + 1
    $r2 = $lastRow + 2

    $ws.Cells.Item($r1, 1).Value = "carbon debt"
    $ws.Cells.Item($r1, 2).Value = "feedstock"
    $ws.Cells.Item($r1, 3).Value = "inflows"
    $ws.Cells.Item($r1, 4).Value = "embodied CO2"
    $ws.Cells.Item($r1, 5).Value = "temp"
    $ws.Cells.Item($r1, 6).Value = "Ratio"
    $ws.Cells.Item($r1, 7).Value = "CO2 removal"

    $ws.Cells.Item($r2, 2).Value = "embodied CO2"
    $ws.Cells.Item($r2, 3).Value = "temp"
    $ws.Cells.Item($r2, 4).Value = "debt CO2"
    $ws.Cells.Item($r2, 5).Value = "outflows"
    $ws.Cells.Item($r2, 6).Value = "Ratio"
    $ws.Cells.Item($r2, 7).Value = "carbon debt factor"

    # mirror the author's navigation: select the newly added block
    $rangeAddr = "A" + $r1 + ":G" + $r2
    $ws.Range($rangeAddr).Select() | Out-Null
}

Add-CarbonDebtRows "Charcoal"
Add-CarbonDebtRows "Syngas"

$wb.Worksheets.Item("Syngas").Activate() | Out-Null
